$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (|S*|/n) across the data rows
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14: Average of SW(S*)/SW(OPT)  -> average of column N
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$b14 = $ws.Range("B14")
$b14.Formula = "=AVERAGE(N2:N11)"
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

# Row 15: Average of SC(S*)/SC(OPT)  -> average of column Z
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$b15 = $ws.Range("B15")
$b15.Formula = "=AVERAGE(Z2:Z11)"

# Row 16: Worst of SW(S*)/SW(OPT) -> min of column N
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$b16 = $ws.Range("B16")
$b16.Formula = "=MIN(N2:N11)"

# Row 17: Worst of SC(S*)/SC(OPT) -> max of column Z
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$b17 = $ws.Range("B17")
$b17.Formula = "=MAX(Z2:Z11)"

# Reuse the exact same formatting (bold, size 12, vertical-centered) from B14
# on B15:B17 without generating duplicate intermediate cell styles.
$b14.Copy()
$b15.PasteSpecial(-4122)
$b16.PasteSpecial(-4122)
$b17.PasteSpecial(-4122)

# Selection now sits on the newly added summary block
$ws.Range("A14:B17").Select()

# Page setup (paper size / orientation) picked up by the resave
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
